# -----------------------------------------------------------------------
# B6-PowerPoint.pptx edit:
#  1) Re-style the three summary tables (slides 14,15,16) from the deck's
#     single custom table style to the built-in
#     "{0DB20450-B7D9-40A0-AAAA-97AE28A8A47F}" table style.
#  2) Swap the slide-master theme's colour scheme from the "Integral /
#     Red Violet" palette over to the stock "Office" palette (the
#     PowerPoint Design theme that was (re)applied).
# -----------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Table style ids -------------------------------------------------
$newTableStyle = "{0DB20450-B7D9-40A0-AAAA-97AE28A8A47F}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2) Theme colour scheme ---------------------------------------------
# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
# RGB values below are packed the way PowerPoint COM expects (0x00BBGGRR).
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$master = $p.Slides.Item(1).Master
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
